# Add a new "msg_count_twitter_engage" stat block between the existing
# "msg_count_twitter" (B:N) and "msg_count_facebook" (O:AA) blocks, and
# refresh the "msg_count_twitter" numbers to the re-computed (engage-excluded)
# values.
#
# Layout before:  B:N = msg_count_twitter | O:AA = msg_count_facebook
# Layout after :  B:N = msg_count_twitter | O:AA = msg_count_twitter_engage (NEW) | AB:AN = msg_count_facebook

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Make room: insert 13 columns at O, shifting the facebook block (O:AA) to AB:AN ---
$ws.Range("O1:AA9").Insert(-4161)  # xlShiftToRight

# --- 2. Label the new block's header (row 1) and merge it, like the other two blocks ---
$ws.Cells.Item(1, 15).Value2 = "msg_count_twitter_engage"

# Merging bordered cells makes Excel synthesize extra split-border styles
# (left/top/bottom/right partial borders) to emulate the merged box visually.
# Clear formatting first so the merge has nothing to split, then re-apply the
# original header style (copied from the sibling block's header) afterwards,
# so every merged cell ends up back on the same style index as before.
$ws.Range("O1:AA1").ClearFormats()
$ws.Range("O1:AA1").Merge()
$ws.Range("B1").Copy()
$ws.Range("O1:AA1").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false

# --- 3. Row 2 sub-headers for the new block (same stat names as the other two blocks) ---
$statLabels = @("sum", "mean", "std", "min", "q25", "median", "q75", "max", "count", "msg_per_mus", "active_mus_n", "active_mus_pc", "active_mus_pc_z")
for ($i = 0; $i -lt $statLabels.Length; $i++) {
    $ws.Cells.Item(2, 15 + $i).Value2 = $statLabels[$i]
}

# --- 4. Refresh msg_count_twitter (B:N) data, rows 4-9, with the recomputed values ---
$twitterData = @{}
$twitterData[4] = @(1479, 61.6, 245.9, 0, 0, 0, 0, 1200, 24, 369.8, 4, 16.7, -1.7)
$twitterData[5] = @(1886800, 753.2, 1772, 0, 0, 125, 881, 43003, 2505, 1191.2, 1584, 63.2, 1.1)
$twitterData[6] = @(50940, 3396, 12233.9, 0, 0, 0, 51.5, 47571, 15, 10188, 5, 33.3, -0.7)
$twitterData[7] = @(51766, 550.7, 1742.2, 0, 0, 0, 525.2, 14888, 94, 1150.4, 45, 47.9, 0.2)
$twitterData[8] = @(315541, 638.7, 2445.9, 0, 0, 2.5, 543.2, 47580, 494, 1247.2, 253, 51.2, 0.4)
$twitterData[9] = @(97303, 470.1, 922, 0, 0, 18, 498.5, 5793, 207, 838.8, 116, 56, 0.7)

foreach ($row in $twitterData.Keys) {
    $vals = $twitterData[$row]
    for ($i = 0; $i -lt $vals.Length; $i++) {
        $ws.Cells.Item($row, 2 + $i).Value2 = $vals[$i]
    }
}

# --- 5. Populate the new msg_count_twitter_engage (O:AA) data, rows 4-9 ---
$engageData = @{}
$engageData[4] = @(500, 20.8, 91.09999999999999, 0, 0, 0, 0, 447, 24, 125, 4, 16.7, -1.6)
$engageData[5] = @(748328, 298.7, 1415.6, 0, 0, 18, 168, 46992, 2505, 464.8, 1610, 64.3, 1.1)
$engageData[6] = @(21724, 1448.3, 5378.6, 0, 0, 0, 9, 20884, 15, 4344.8, 5, 33.3, -0.7)
$engageData[7] = @(18156, 193.1, 1081.3, 0, 0, 0, 64.8, 10421, 94, 422.2, 43, 45.7, 0.1)
$engageData[8] = @(101125, 204.7, 868.9, 0, 0, 2, 92.2, 9692, 494, 392, 258, 52.2, 0.4)
$engageData[9] = @(19229, 92.90000000000001, 281.9, 0, 0, 3, 68.5, 3183, 207, 164.4, 117, 56.5, 0.7)

foreach ($row in $engageData.Keys) {
    $vals = $engageData[$row]
    for ($i = 0; $i -lt $vals.Length; $i++) {
        $ws.Cells.Item($row, 15 + $i).Value2 = $vals[$i]
    }
}

Write-Output "edit applied"
